$wb = $excel.ActiveWorkbook

# --- Overview sheet: update row 3 (73e28401 file) with new handoff status/date ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-16 02:43:42"

# --- zh-cn sheet: update row 3 (73e28401 file) Status/Handoff-Datetime/Error Detail, widen column P ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-16 02:43:37"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b02ae85d0352b4815ad8b707d677cda5611c8474/e2e/73e28401-ac3f-4dda-8550-b29fa5410a52.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc1f57c5dacdd9c4ce7d9346ff87164dbc4fe9c4/e2e/73e28401-ac3f-4dda-8550-b29fa5410a52.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666

# --- de-de sheet: update row 3 (73e28401 file) Status/Handoff-Datetime/Error Detail, widen column P ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-16 02:43:42"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b02ae85d0352b4815ad8b707d677cda5611c8474/e2e/73e28401-ac3f-4dda-8550-b29fa5410a52.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc1f57c5dacdd9c4ce7d9346ff87164dbc4fe9c4/e2e/73e28401-ac3f-4dda-8550-b29fa5410a52.md."
$dede.Columns.Item(16).ColumnWidth = 39.1666666
